# "updated after lecture 4"
# Slide 2's title placeholder reads "Power BI Quick Start #4 " + "E01".
# Renumber it to "Power BI Quick Start #2 " + "E04", editing each existing
# run's text individually so the run-level formatting (rPr) is preserved
# exactly as authored.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$title = $s.Shapes.Title
$tr = $title.TextFrame.TextRange

$tr.Runs(1).Text = "Power BI Quick Start #2 "
$tr.Runs(2).Text = "E04"
